$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-08-26T17:07:33"

$ws.Range("U4").Value = 55
$ws.Range("V4").Value = 45.3
$ws.Range("W4").Value = 35.59
$ws.Range("X4").Value = 32.66
$ws.Range("Y4").Value = 29.11
$ws.Range("Z4").Value = 28.86

$ws.Range("U6").Value = -2.25
$ws.Range("V6").Value = -1.9
$ws.Range("W6").Value = -1.42
$ws.Range("X6").Value = -0.98
$ws.Range("Y6").Value = -0.84
$ws.Range("Z6").Value = -0.84

$ws.Range("U9").Value = 55.91
$ws.Range("V9").Value = 50.37
$ws.Range("W9").Value = 36.15
$ws.Range("X9").Value = 33.95
$ws.Range("Y9").Value = 30.6
$ws.Range("Z9").Value = 30.43

$ws.Range("U11").Value = -1.34
$ws.Range("V11").Value = -0.28
$ws.Range("W11").Value = -0.87
$ws.Range("X11").Value = 0.31
$ws.Range("Y11").Value = 0.64
$ws.Range("Z11").Value = 0.73

$ws.Range("V12").Value = 3.45

$ws.Range("U14").Value = 55.91
$ws.Range("V14").Value = 50.37
$ws.Range("W14").Value = 36.15
$ws.Range("Y14").Value = 30.6
$ws.Range("Z14").Value = 30.43

$ws.Range("U16").Value = -1.34
$ws.Range("V16").Value = -0.28
$ws.Range("W16").Value = -0.87
$ws.Range("Z16").Value = 0.73

$ws.Range("V17").Value = 3.45

$ws.Range("U19").Value = 54.74
$ws.Range("V19").Value = 45.17
$ws.Range("W19").Value = 35.59
$ws.Range("Y19").Value = 29.22
$ws.Range("Z19").Value = 28.94

$ws.Range("U21").Value = -2.52
$ws.Range("V21").Value = -2.03
$ws.Range("W21").Value = -1.42
$ws.Range("Y21").Value = -0.73
$ws.Range("Z21").Value = -0.75

$ws.Range("U24").Value = 54.74
$ws.Range("V24").Value = 45.17
$ws.Range("W24").Value = 35.59
$ws.Range("Y24").Value = 29.22
$ws.Range("Z24").Value = 28.94

$ws.Range("U26").Value = -2.52
$ws.Range("V26").Value = -2.03
$ws.Range("W26").Value = -1.42
$ws.Range("Y26").Value = -0.73
$ws.Range("Z26").Value = -0.75

$ws.Range("U29").Value = 54.42
$ws.Range("V29").Value = 45
$ws.Range("W29").Value = 35.56
$ws.Range("X29").Value = 32.85
$ws.Range("Y29").Value = 29.31
$ws.Range("Z29").Value = 29.03

$ws.Range("U31").Value = -2.83
$ws.Range("V31").Value = -2.2
$ws.Range("W31").Value = -1.46
$ws.Range("X31").Value = -0.79
$ws.Range("Y31").Value = -0.64
$ws.Range("Z31").Value = -0.67

$ws.Range("U34").Value = 57.37
$ws.Range("V34").Value = 48.37
$ws.Range("W34").Value = 36.54
$ws.Range("X34").Value = 34.86
$ws.Range("Y34").Value = 31.56
$ws.Range("Z34").Value = 31.42

$ws.Range("U36").Value = 0.11
$ws.Range("V36").Value = 1.16
$ws.Range("X36").Value = 1.22

$ws.Range("U39").Value = 55
$ws.Range("V39").Value = 45.3
$ws.Range("W39").Value = 35.59
$ws.Range("X39").Value = 32.66
$ws.Range("Y39").Value = 29.11
$ws.Range("Z39").Value = 28.86

$ws.Range("U41").Value = -2.25
$ws.Range("V41").Value = -1.9
$ws.Range("W41").Value = -1.42
$ws.Range("X41").Value = -0.98
$ws.Range("Y41").Value = -0.84
$ws.Range("Z41").Value = -0.84

$ws.Range("U44").Value = 57.25
$ws.Range("V44").Value = 47.11
$ws.Range("W44").Value = 36.9
$ws.Range("Y44").Value = 29.89
$ws.Range("Z44").Value = 29.7

$ws.Range("V46").Value = -0.09
$ws.Range("W46").Value = -0.11
$ws.Range("Y46").Value = -0.06

$ws.Range("U49").Value = 59.7
$ws.Range("V49").Value = 49.02
$ws.Range("W49").Value = 38.28
$ws.Range("Y49").Value = 30.13
$ws.Range("Z49").Value = 29.97

$ws.Range("U51").Value = 2.45
$ws.Range("V51").Value = 1.81
$ws.Range("W51").Value = 1.26

$ws.Range("U54").Value = 56.8
$ws.Range("V54").Value = 47.21
$ws.Range("W54").Value = 37.16
$ws.Range("X54").Value = 34.08
$ws.Range("Y54").Value = 30.97
$ws.Range("Z54").Value = 30.9

$ws.Range("U56").Value = -0.45
$ws.Range("W56").Value = 0.15
$ws.Range("X56").Value = 0.44
$ws.Range("Y56").Value = 1.02
$ws.Range("Z56").Value = 1.21

$ws.Range("U59").Value = 59.33
$ws.Range("V59").Value = 48.82
$ws.Range("W59").Value = 38.24
$ws.Range("Y59").Value = 30.88
$ws.Range("Z59").Value = 30.74

$ws.Range("U61").Value = 2.08
$ws.Range("V61").Value = 1.61
$ws.Range("W61").Value = 1.22
$ws.Range("Y61").Value = 0.93
$ws.Range("Z61").Value = 1.05

$ws.Range("U64").Value = 60.33
$ws.Range("V64").Value = 49.64
$ws.Range("W64").Value = 38.84
$ws.Range("X64").Value = 35.26
$ws.Range("Y64").Value = 31.3
$ws.Range("Z64").Value = 31.16

$ws.Range("U66").Value = 3.08
$ws.Range("V66").Value = 2.43
$ws.Range("W66").Value = 1.83
$ws.Range("X66").Value = 1.62
$ws.Range("Y66").Value = 1.35
$ws.Range("Z66").Value = 1.46

$ws.Range("U69").Value = 61.04
$ws.Range("V69").Value = 50.22
$ws.Range("W69").Value = 39.29
$ws.Range("Y69").Value = 31.73
$ws.Range("Z69").Value = 31.59

$ws.Range("U71").Value = 3.78
$ws.Range("V71").Value = 3.01
$ws.Range("W71").Value = 2.28

$ws.Range("U74").Value = 59.21
$ws.Range("V74").Value = 48.87
$ws.Range("W74").Value = 38.24
$ws.Range("Y74").Value = 30.75
$ws.Range("Z74").Value = 30.61

$ws.Range("U76").Value = 1.95
$ws.Range("V76").Value = 1.66
$ws.Range("W76").Value = 1.22

$ws.Range("U79").Value = 59.6
$ws.Range("V79").Value = 49.12
$ws.Range("W79").Value = 38.46
$ws.Range("Y79").Value = 30.97
$ws.Range("Z79").Value = 30.8

$ws.Range("U81").Value = 2.35
$ws.Range("V81").Value = 1.92
$ws.Range("W81").Value = 1.45
$ws.Range("Y81").Value = 1.01

$ws.Range("U84").Value = 57.83
$ws.Range("V84").Value = 48.32
$ws.Range("W84").Value = 38.12
$ws.Range("X84").Value = 34.93
$ws.Range("Y84").Value = 31.56
$ws.Range("Z84").Value = 31.29

$ws.Range("U86").Value = 0.58
$ws.Range("V86").Value = 1.11
$ws.Range("W86").Value = 1.11
$ws.Range("X86").Value = 1.29
$ws.Range("Z86").Value = 1.6

$ws.Range("U89").Value = 54.42
$ws.Range("V89").Value = 44.96
$ws.Range("W89").Value = 35.56
$ws.Range("X89").Value = 32.85
$ws.Range("Y89").Value = 29.31
$ws.Range("Z89").Value = 29.03

$ws.Range("U91").Value = -2.83
$ws.Range("V91").Value = -2.25
$ws.Range("W91").Value = -1.46
$ws.Range("X91").Value = -0.79
$ws.Range("Y91").Value = -0.64
$ws.Range("Z91").Value = -0.67
